# Feature: update Excel export formats for pre-inscriptions and references
#
# The "Listado de Referencias Pendientes" export template shipped with two
# sample/demo rows of data underneath the table header (row 9: "Daniel
# Wolke..."; row 10: "jossue irias..."). Those were only placeholder/test
# rows used while building the export format and must not ship in the
# template - the table should start empty, ready for the real export data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the sample/demo data rows (keep header row 8 + row formatting/style
# intact) so the exported table template starts with blank rows, same as
# rows 11:16 below them.
$ws.Range("B9:K10").ClearContents()

# Restore the normal (non-demo) selection/active cell that Excel leaves in
# the file on save.
$ws.Range("J24").Select()
